# Localization status report refresh:
#  - Status moves from "Ready for handoff" to "In Translation" for the
#    tracked file (zh-cn / de-de), reflected on all three sheets.
#  - The narrower status text lets the status columns shrink, so the
#    columns are resized to match the new (auto-fitted) content width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 13.4101845877511

# --- Overview sheet: zh-cn / de-de status columns (E & F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns("E:E").ColumnWidth = $newWidth
$overview.Columns("F:F").ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns("C:C").ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns("C:C").ColumnWidth = $newWidth
